$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a "last updated" date serial that was bumped
# forward by one day (45180 -> 45181) for every data row (rows 2 through 554).
for ($r = 2; $r -le 554; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value = 45181
    }
}
